# The authored change swaps the two theme parts of the deck:
#   ppt/theme/theme1.xml  "Office Theme" (Office colours)  <->  ppt/theme/theme2.xml "Integral" (green/olive colours)
#
# ppt/theme/theme2.xml is the theme actually driving every slide (it is the
# presentation's primary theme, wired up from slideMaster1 / presentation.xml),
# so the visible effect of the swap is that the deck's theme colours change
# from the "Integral" palette to the standard "Office" palette. Apply that
# through the PowerPoint colour-scheme object model, which writes straight
# back into the clrScheme of the theme part backing the slide master.

$p = $ppt.ActivePresentation

# Target palette == the "Office" colour scheme that currently lives in
# ppt/theme/theme1.xml, expressed as RRGGBB, in clrScheme order.
$targetColors = @(
  @{ Index = 1;  Name = "dk1";      Hex = "000000" }
  @{ Index = 2;  Name = "lt1";      Hex = "FFFFFF" }
  @{ Index = 3;  Name = "dk2";      Hex = "44546A" }
  @{ Index = 4;  Name = "lt2";      Hex = "E7E6E6" }
  @{ Index = 5;  Name = "accent1";  Hex = "5B9BD5" }
  @{ Index = 6;  Name = "accent2";  Hex = "ED7D31" }
  @{ Index = 7;  Name = "accent3";  Hex = "A5A5A5" }
  @{ Index = 8;  Name = "accent4";  Hex = "FFC000" }
  @{ Index = 9;  Name = "accent5";  Hex = "4472C4" }
  @{ Index = 10; Name = "accent6";  Hex = "70AD47" }
  @{ Index = 11; Name = "hlink";    Hex = "0563C1" }
  @{ Index = 12; Name = "folHlink"; Hex = "954F72" }
)

# The theme colour scheme, reached from the slide master so the write lands
# on the theme part (not a one-off per-slide override).
$themeColors = $p.SlideMaster.Theme.ThemeColorScheme

foreach ($c in $targetColors) {
    $hex = $c.Hex
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    # PowerPoint RGB() packs as 0x00BBGGRR.
    $comRgb = ($b * 65536) + ($g * 256) + $r

    $themeColors.Colors($c.Index).RGB = $comRgb
    Write-Host ("Set {0} (#{1}) -> {2}" -f $c.Name, $c.Index, $hex)
}

Write-Host "Theme colour scheme updated to the Office palette."
